# The sheet used to have an extra header row ("Number" / "String") above the
# data rows, which forced callers to declare a column type. Remove that row
# so the "ID" column starts immediately below the "ID"/"Name"/"Location"
# header row - Excel reflows the remaining rows upward and drops the now
# superfluous trailing blank row from the used range automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
